$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 661; this shifts the existing rows 661-714 down to 662-715.
$ws.Rows.Item(661).Insert()

# Populate the newly inserted row 661 with the new data record.
$ws.Range("A661").Value = 3
$ws.Range("B661").Value = "Femacal de La Calera"
$ws.Range("C661").Value = "Coquimbo"
$ws.Range("D661").Value = 45265
$ws.Range("E661").Value = 5
$ws.Range("F661").Value = 100112031
$ws.Range("G661").Value = "Poroto verde"
$ws.Range("H661").Value = "Magnum"
$ws.Range("I661").Value = "Primera"
$ws.Range("J661").Value = 40
$ws.Range("K661").Value = 35000
$ws.Range("L661").Value = 35000
$ws.Range("M661").Value = 35000
$ws.Range("N661").Value = "`$/malla 25 kilos"
$ws.Range("O661").Value = "Provincia de Limarí"
$ws.Range("P661").Value = 1400
$ws.Range("Q661").Value = 25
$ws.Range("R661").Value = "Hortaliza"
